$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data points (refreshed source values)
$ws.Range("B205").Value = 4373594970000
$ws.Range("D205").Value = 195008452486.4765
$ws.Range("B206").Value = 4411934620000
$ws.Range("D206").Value = 196098882059.7977
$ws.Range("B210").Value = 4577407590000
$ws.Range("D210").Value = 197829026890.5965
$ws.Range("B211").Value = 4566459490000
$ws.Range("D211").Value = 204213160361.4111
$ws.Range("B212").Value = 4592275590000
$ws.Range("D212").Value = 198028441162.9506
$ws.Range("B213").Value = 4639859400000
$ws.Range("D213").Value = 196360924007.7125
$ws.Range("B214").Value = 4641345140000
$ws.Range("D214").Value = 196115881881.1865
$ws.Range("B215").Value = 4681223420000
$ws.Range("D215").Value = 204604827761.64
$ws.Range("B216").Value = 4725508480000
$ws.Range("D216").Value = 203831027810.3441
$ws.Range("B217").Value = 4680322510000
$ws.Range("D217").Value = 206432055012.9981
$ws.Range("B219").Value = 4809150480000
$ws.Range("D219").Value = 208310072120.0702
$ws.Range("B221").Value = 4958595660000
$ws.Range("D221").Value = 200173330572.0967
$ws.Range("B222").Value = 5004666910000
$ws.Range("D222").Value = 206283339234.405
$ws.Range("B223").Value = 5020790900000
$ws.Range("D223").Value = 211475627271.8466
$ws.Range("B224").Value = 5059232680000
$ws.Range("D224").Value = 227082617645.2899
$ws.Range("B225").Value = 5094308060000
$ws.Range("D225").Value = 231533988533.9885
$ws.Range("B226").Value = 5178041490000
$ws.Range("D226").Value = 224145123064.5561
$ws.Range("B227").Value = 5214187690000
$ws.Range("D227").Value = 222760943901.6367
$ws.Range("B228").Value = 5235568230000
$ws.Range("D228").Value = 237994978328.2156
$ws.Range("B229").Value = 5179738620000
$ws.Range("D229").Value = 241084591915.8116
$ws.Range("B230").Value = 5290478980000
$ws.Range("D230").Value = 246567659217.4866
$ws.Range("B231").Value = 5390398340000
$ws.Range("D231").Value = 248749920512.0455
$ws.Range("B232").Value = 5449356120000
$ws.Range("D232").Value = 244698227636.6543
$ws.Range("B233").Value = 5471474170000
$ws.Range("D233").Value = 254209314005.6218
$ws.Range("B234").Value = 5507491430000
$ws.Range("D234").Value = 264768603178.1897
$ws.Range("B236").Value = 5564521500000
$ws.Range("D236").Value = 259066790507.9822
$ws.Range("B239").Value = 5617130550000
$ws.Range("D239").Value = 253143058331.5119
$ws.Range("B240").Value = 5647837280000
$ws.Range("D240").Value = 250967585233.9003
$ws.Range("B241").Value = 5542014840000
$ws.Range("D241").Value = 253256051601.582
$ws.Range("B242").Value = 5630383690000
$ws.Range("D242").Value = 259884453460.5434
$ws.Range("B243").Value = 5704249840000
$ws.Range("D243").Value = 254295839565.2552
$ws.Range("B244").Value = 5739159050000
$ws.Range("D244").Value = 260129496843.0881
$ws.Range("B245").Value = 5742427260000
$ws.Range("D245").Value = 246010541422.8308
$ws.Range("B246").Value = 5825723830000
$ws.Range("D246").Value = 252980514452.8375
$ws.Range("B247").Value = 5801917230000
$ws.Range("D247").Value = 245858014335.5658
$ws.Range("B248").Value = 5833040250000
$ws.Range("D248").Value = 242345086968.173
$ws.Range("B249").Value = 5855415460000
$ws.Range("D249").Value = 240184235925.8127
$ws.Range("B250").Value = 5887405600000
$ws.Range("D250").Value = 234668314187.5463
$ws.Range("B251").Value = 5915934540000
$ws.Range("D251").Value = 238841410132.7046
$ws.Range("B252").Value = 5940210650000
$ws.Range("D252").Value = 253955001539.0667
$ws.Range("B256").Value = 6077524080000
$ws.Range("D256").Value = 280799406754.406
$ws.Range("B258").Value = 6224248910000
$ws.Range("D258").Value = 280487810714.3501

# Append new row 259 (2023-06-01 data point)
$ws.Range("A258").Copy($ws.Range("A259"))
$ws.Range("A259").Value = 45078
$ws.Range("B259").Value = 6243183470000
$ws.Range("C259").Value = 0.04592443593311565
$ws.Range("D259").Value = 286714679286.7017
